$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New quote row (row 32): Trump tweet about the cure being worse than the
# problem, added to the coronavirus quotes tracker.
$ws.Range("A32").Value = 43912
$ws.Range("A32").NumberFormat = $ws.Range("A2").NumberFormat

$ws.Range("B32").Value = "WE CANNOT LET THE CURE BE WORSE THAN THE PROBLEM ITSELF."

$ws.Range("C32").Value = "Donald Trump"

$ws.Range("D32").Value = "https://twitter.com/realDonaldTrump/status/1241935285916782593"
$ws.Range("D32").Style = $ws.Range("D2").Style
$ws.Hyperlinks.Add($ws.Range("D32"), "https://twitter.com/realDonaldTrump/status/1241935285916782593")

$ws.Range("E32").Value = "Quote"

$ws.Range("G32").Value = 1

# Move the live selection down to the newly added row, matching the
# author's cursor position after the edit.
$ws.Range("E33").Select()
